# edit.ps1 - Applies the "Added mini quests to list" change to AssetList.docx
#
# Summary of changes:
#  1. Split the "Sharky the Goldfish" bullet into two bullets: a new
#     "Lily the Lilac" bullet followed by a cleaned-up "Sharky the Goldfish"
#     bullet (single run, no proofErr spell-check markers).
#  2. Clean up the "Sharky's Lost uPed" bullet into a single run (no
#     proofErr spell-check markers).
#  3. Change "Tulip fertilizer" bullet text to "Lilac fertilizer" (kept as
#     two runs: "Lilac" + " fertilizer").
#  4. Insert a new "Remove item from inventory (quest)" bullet under "Key
#     items" (right after "Toggle key item inventory on/off", before
#     "Character status"), and move the _GoBack bookmark there (it used to
#     sit after the very last bullet, "Exit").
#  5. Clean up the "Enemy aggro/de-aggro" bullet into a single run (no
#     proofErr spell-check markers).

$d = $word.ActiveDocument

function Set-ParaXml($para, $innerBodyXml) {
    # Replaces the contents of a paragraph's range with a precise OOXML
    # fragment, bypassing Find/Replace's tendency to leave stray markup
    # (e.g. orphaned <w:proofErr/> elements) behind.
    $full = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($full)
}

# ---------------------------------------------------------------------
# 1. "Sharky the Goldfish" -> "Lily the Lilac" (new) + "Sharky the Goldfish"
# ---------------------------------------------------------------------
$sharkyFish = $d.Content.Find.Execute("Sharky the Goldfish", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pSharkyFish = $d.Paragraphs(16)

$pSharkyFish.Range.InsertParagraphBefore()
$lily = $d.Paragraphs(16)
$lily.Range.Text = "Lily the Lilac"

$sharkyFish2 = $d.Paragraphs(17)
Set-ParaXml $sharkyFish2 '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Sharky the Goldfish</w:t></w:r></w:p></w:body>'

# ---------------------------------------------------------------------
# 2. "Sharky's Lost uPed" -> single run, no proofErr
# ---------------------------------------------------------------------
$pLostUped = $d.Paragraphs(69)
Set-ParaXml $pLostUped '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Sharky&#8217;s Lost uPed</w:t></w:r></w:p></w:body>'

# ---------------------------------------------------------------------
# 3. "Tulip fertilizer" -> "Lilac" + " fertilizer" (two runs)
# ---------------------------------------------------------------------
$pFertilizer = $d.Paragraphs(70)
Set-ParaXml $pFertilizer '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Lilac</w:t></w:r><w:r><w:t xml:space="preserve"> fertilizer</w:t></w:r></w:p></w:body>'

# ---------------------------------------------------------------------
# 4. Insert "Remove item from inventory (quest)" after "Toggle key item
#    inventory on/off" and before "Character status"; move the _GoBack
#    bookmark from the "Exit" paragraph to the new paragraph.
# ---------------------------------------------------------------------
$pToggleKeyItem = $d.Paragraphs(103)
if ($pToggleKeyItem.Range.Text -notlike "Toggle key item inventory on/off*") {
    throw "Unexpected paragraph located where 'Toggle key item inventory on/off' was expected: $($pToggleKeyItem.Range.Text)"
}

# Insert after "Toggle key item inventory on/off" (rather than before
# "Character status") so the new bullet inherits the ilvl=1 indent level.
$pToggleKeyItem.Range.InsertParagraphAfter()
$pRemoveItem = $d.Paragraphs(104)
$pRemoveItem.Range.Text = "Remove item from inventory (quest)"

# Move the _GoBack bookmark onto the new paragraph (it currently trails
# the very last bullet, "Exit"). Exclude the trailing paragraph mark from
# the bookmark range so both bookmarkStart and bookmarkEnd land inside
# this paragraph, right after the run (matching Word's normal behavior).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($pRemoveItem.Range.Start, $pRemoveItem.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 5. "Enemy aggro/de-aggro" -> single run, no proofErr
# ---------------------------------------------------------------------
$pAggro = $d.Paragraphs(113)
if ($pAggro.Range.Text -notlike "Enemy aggro*") {
    throw "Unexpected paragraph located where 'Enemy aggro/de-aggro' was expected: $($pAggro.Range.Text)"
}
Set-ParaXml $pAggro '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Enemy aggro/de-aggro</w:t></w:r></w:p></w:body>'
